$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header cells for the new columns
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the style of an existing header cell (e.g. AC1) so the new headers match formatting
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Fill in the team record values for every data row (2-53)
for ($r = 2; $r -le 53; $r++) {
    $ws.Cells.Item($r, 30).Value = 65   # AD - Wins
    $ws.Cells.Item($r, 31).Value = 97   # AE - Losses
    $ws.Cells.Item($r, 32).Value = 0    # AF - Ties
}
